# Append new ticker rows to the worksheet, extending the used range
# from A1:A176 to A1:A179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A177").Value = "IMX-USD"
$ws.Range("A178").Value = "TAO-USD"
$ws.Range("A179").Value = "GRT-USD"
